$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Rewrite the "cannons" feedback-response paragraph.
#    Old: three runs describing non-normalised vector movement.
#    New: four sentences describing matrix-based parent/child movement.
# ------------------------------------------------------------------
$cannonsRange = $d.Content
$found = $cannonsRange.Find.Execute("As for the cannons moving out of sync with the boss, implemented movement using non-normalised vectors, and made added a logic check to make the cannons know they are in position before moving left and right which stops restarting breaking it.")
if ($found) {
    $cannonsRange.Text = "As for the cannons moving out of sync with the boss, I have implemented movement using matrices from my library. I have made the boss body the parent and the turrets the children preventing unwanted movement."
}

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the document (last
#    paragraph, after the "Fixed the score issue..." text) up to the
#    end of the "Perhaps a space invaders deal..." paragraph, which is
#    now the most-recently-edited spot in the document.
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Perhaps a space invaders deal, wait for all enemies in the array are turned off, before turning them all on again.") | Out-Null
$target.Collapse(0)

# A zero-length Range placed exactly on a paragraph mark cannot be
# bookmarked directly, so insert a temporary marker character, wrap the
# bookmark around it, then delete the marker -- leaving the bookmark
# collapsed in the correct spot, right after the run and before </w:p>.
$target.InsertAfter("@")
$markerRange = $d.Range($target.Start, $target.Start + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange = $d.Range($target.Start, $target.Start + 1)
$markerRange.Text = ""
